# Update "want to go" counts (column F) across sheets to reflect the
# regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 20
$wsExhibit.Range("F8").Value  = 7856
$wsExhibit.Range("F9").Value  = 749
$wsExhibit.Range("F10").Value = 219
$wsExhibit.Range("F11").Value = 1093
$wsExhibit.Range("F12").Value = 731
$wsExhibit.Range("F15").Value = 193
$wsExhibit.Range("F18").Value = 810

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 1

# --- Sheet "全部类型" (All types, combined view) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 20
$wsAll.Range("F7").Value  = 1
$wsAll.Range("F9").Value  = 7856
$wsAll.Range("F10").Value = 749
$wsAll.Range("F11").Value = 219
$wsAll.Range("F12").Value = 1093
$wsAll.Range("F13").Value = 731
$wsAll.Range("F16").Value = 193
$wsAll.Range("F19").Value = 810
